$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells keep their original text formatting
# (values like "350.70" would otherwise be auto-converted to numbers by Excel).

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '50.868.20'
$ws.Range('E2').Value = '  -2.12%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.750.18'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '350.70'
$ws.Range('E5').Value = '  -2.27%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '107.04'
$ws.Range('E6').Value = '  -2.52%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.544'
$ws.Range('E7').Value = '  -2.99%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.577'
$ws.Range('E9').Value = '  -2.45%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '38.98'
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('E11').Value = '  +3.37%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0829'
$ws.Range('E12').Value = '  -3.08%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.65'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('E14').Value = '  -2.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.181.82'
$ws.Range('E15').Value = '  -0.90%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.749.10'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '50.817.29'
$ws.Range('E18').Value = '  -1.88%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.53'
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('E20').Value = '  -3.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.93'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('E22').Value = '  -3.27%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '69.12'
$ws.Range('E23').Value = '  -0.92%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '263.11'
$ws.Range('E24').Value = '  -3.92%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.69'
$ws.Range('E25').Value = '  -2.46%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '25.76'
$ws.Range('E27').Value = '  -3.28%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.161'
$ws.Range('E28').Value = '  +13.46%  '
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '51.41'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.99'
$ws.Range('E32').Value = '  +4.00%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '34.20'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.0440'
$ws.Range('E34').Value = '  -6.82%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '5.30'
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0826'
$ws.Range('E36').Value = '  -1.55%  '
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '18.21'
$ws.Range('E38').Value = '  +0.29%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.12'
$ws.Range('E39').Value = '  -2.69%  '
$ws.Range('E40').Value = '  -3.36%  '
$ws.Range('E41').Value = '  -1.67%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.46'
$ws.Range('E42').Value = '  -6.11%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '120.39'
$ws.Range('E43').Value = '  -3.73%  '
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '21.87'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.076.77'
$ws.Range('E46').Value = '  +0.56%  '
$ws.Range('E47').Value = '  -1.30%  '
$ws.Range('E48').Value = '  -1.61%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.904'
$ws.Range('E49').Value = '  -4.72%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '5.40'
$ws.Range('E50').Value = '  -5.62%  '
$ws.Range('E51').Value = '  +4.36%  '
